# Apply the batch_demo.xlsx edit: duplicate row 4 into rows 5-7 (io_series with
# multiple "*" entries), then adjust the C/D descriptor columns and the
# numeric "E" parameter column so each row represents a different value,
# and fix up row 3's values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate row 4 into rows 5, 6 and 7 (copy formatting + values) ---
$srcRow = $ws.Range("A4:L4")
$srcRow.Copy()
$ws.Range("A5:L5").PasteSpecial()
$srcRow.Copy()
$ws.Range("A6:L6").PasteSpecial()
$srcRow.Copy()
$ws.Range("A7:L7").PasteSpecial()
$excel.CutCopyMode = 0

# --- Row 3: switch descriptor from "comparison"/"med" to "blurred"/"reg" ---
$ws.Range("C3").Value = "blurred"
$ws.Range("D3").Value = "reg"
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 3

# --- Row 4: correct the size/param values ---
$ws.Range("B4").Value = 51
$ws.Range("E4").Value = 15

# --- Row 5 ---
$ws.Range("B5").Value = 51
$ws.Range("C5").Value = "blurred"
$ws.Range("D5").Value = "reg"
$ws.Range("E5").Value = 20

# --- Row 6 ---
$ws.Range("B6").Value = 51
$ws.Range("C6").Value = "blurred"
$ws.Range("D6").Value = "reg"
$ws.Range("E6").Value = 25

# --- Row 7 ---
$ws.Range("B7").Value = 51
$ws.Range("C7").Value = "blurred"
$ws.Range("D7").Value = "reg"
$ws.Range("E7").Value = 30

# --- Update selection to match the final cursor position ---
$ws.Range("E7").Select() | Out-Null
